# Update weekly Fruta/hortaliza price data for "Terminal La Palmera de La Serena - Coco"
# (shuffles Fecha, Volumen, Precio minimo/maximo/promedio ponderado, and Precio $/Kg
# values across the data rows, per upstream weekly data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44333
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 19500
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19750
$ws.Range("S2").Value = 988

$ws.Range("D3").Value = 44350
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("S3").Value = 975

$ws.Range("D4").Value = 44335
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19500
$ws.Range("S4").Value = 975

$ws.Range("D5").Value = 44474
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19500
$ws.Range("S5").Value = 975

$ws.Range("D6").Value = 44427
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20500
$ws.Range("S6").Value = 1025

$ws.Range("D7").Value = 44407
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 20500
$ws.Range("S7").Value = 1025

$ws.Range("D8").Value = 44420
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("S8").Value = 1025

$ws.Range("D9").Value = 44434
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("S9").Value = 1025

$ws.Range("D10").Value = 44441
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 21000
$ws.Range("P10").Value = 20500
$ws.Range("S10").Value = 1025

$ws.Range("D11").Value = 44336
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 19500
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 19750
$ws.Range("S11").Value = 988

$ws.Range("D12").Value = 44442
$ws.Range("M12").Value = 140
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 21000
$ws.Range("P12").Value = 20500
$ws.Range("S12").Value = 1025

$ws.Range("D13").Value = 44467
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20500
$ws.Range("S13").Value = 1025

$ws.Range("D14").Value = 44301
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 19000
$ws.Range("P14").Value = 18500
$ws.Range("S14").Value = 925

$ws.Range("D15").Value = 44410
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 20500
$ws.Range("S15").Value = 1025

$ws.Range("D16").Value = 44435
$ws.Range("M16").Value = 260
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 22000
$ws.Range("P16").Value = 21115
$ws.Range("S16").Value = 1056

$ws.Range("D17").Value = 44364
$ws.Range("M17").Value = 140
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 21000
$ws.Range("P17").Value = 20500
$ws.Range("S17").Value = 1025

$ws.Range("D18").Value = 44417
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 21000
$ws.Range("P18").Value = 20500
$ws.Range("S18").Value = 1025

$ws.Range("D19").Value = 44365
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 21000
$ws.Range("P19").Value = 20500
$ws.Range("S19").Value = 1025

$ws.Range("D20").Value = 44462
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 19500
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 19750
$ws.Range("S20").Value = 988

$ws.Range("D21").Value = 44431
$ws.Range("M21").Value = 160
$ws.Range("N21").Value = 21000
$ws.Range("O21").Value = 22000
$ws.Range("P21").Value = 21500
$ws.Range("S21").Value = 1075

$ws.Range("D22").Value = 44418
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 20000
$ws.Range("O22").Value = 21000
$ws.Range("P22").Value = 20500
$ws.Range("S22").Value = 1025

$ws.Range("D23").Value = 44466
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 20000
$ws.Range("O23").Value = 21000
$ws.Range("P23").Value = 20500
$ws.Range("S23").Value = 1025

$ws.Range("D24").Value = 44445
$ws.Range("M24").Value = 160
$ws.Range("N24").Value = 20000
$ws.Range("O24").Value = 21000
$ws.Range("P24").Value = 20500
$ws.Range("S24").Value = 1025

$ws.Range("D25").Value = 44428
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 20000
$ws.Range("O25").Value = 21000
$ws.Range("P25").Value = 20500
$ws.Range("S25").Value = 1025

$ws.Range("D26").Value = 44448
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 20000
$ws.Range("O26").Value = 21000
$ws.Range("P26").Value = 20500
$ws.Range("S26").Value = 1025

$ws.Range("D27").Value = 44315
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 20000
$ws.Range("O27").Value = 21000
$ws.Range("P27").Value = 20500
$ws.Range("S27").Value = 1025

$ws.Range("D28").Value = 44473
$ws.Range("M28").Value = 40
$ws.Range("N28").Value = 19500
$ws.Range("O28").Value = 20000
$ws.Range("P28").Value = 19750
$ws.Range("S28").Value = 988

$ws.Range("D29").Value = 44326
$ws.Range("M29").Value = 160
$ws.Range("N29").Value = 19500
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 19750
$ws.Range("S29").Value = 988

$ws.Range("D30").Value = 44343
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 19500
$ws.Range("O30").Value = 20000
$ws.Range("P30").Value = 19750
$ws.Range("S30").Value = 988
